# Ticket INC-00234.docx: correct the accented "responsável" JSON key to
# the plain-ASCII "responsavel" spelling used consistently elsewhere in
# this JSON-ish snippet (it now matches "emailResponsavel"). Touching
# that character mid-document nudges Word's internal "last edit position"
# (_GoBack) bookmark down to the e-mail line below, and the run there
# gets re-split as well, so both paragraphs end up re-serialized with a
# few more runs and fresh spell-check (proofErr) markers.

$d = $word.ActiveDocument

function Set-RangeXml($range, [string]$bodyInnerXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?>' +
           '<?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyInnerXml + '</w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# --- Paragraph 1 ("numeroIncidente") --------------------------------------
# Text is unchanged, but this paragraph loses the "_GoBack" bookmark - it
# relocates into the e-mail paragraph below (a document can only have one
# bookmark with a given name).
$p1 = $d.Paragraphs.Item(1)
$p1xml = '<w:p>' +
         '<w:r><w:t>"</w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>numeroIncidente</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t>": "00234",</w:t></w:r>' +
         '</w:p>'
Set-RangeXml $p1.Range $p1xml

# --- Paragraphs 4 + 5 ("responsável" -> "responsavel", "emailResponsavel")
# Rewritten together in a single InsertXML call: paragraph 5 is the last
# paragraph in the body, and replacing only its own range there leaves a
# stray empty trailing paragraph behind, so the safe way to touch it is to
# span both paragraphs 4 and 5 in one go.
$p4 = $d.Paragraphs.Item(4)
$p5 = $d.Paragraphs.Item(5)
$combined = $d.Range($p4.Range.Start, $p5.Range.End)

$p45xml = '<w:p>' +
          '<w:r><w:t>"</w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:t>respons</w:t></w:r>' +
          '<w:r><w:t>a</w:t></w:r>' +
          '<w:r><w:t>vel</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t>": "Caio Henrique",</w:t></w:r>' +
          '</w:p>' +
          '<w:p>' +
          '<w:r><w:t>"</w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:t>emailR</w:t></w:r>' +
          '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
          '<w:bookmarkEnd w:id="0"/>' +
          '<w:r><w:t>esponsavel</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t>": "</w:t></w:r>' +
          '<w:r><w:t>caio.</w:t></w:r>' +
          '<w:r><w:t>henrique@gmail.com"</w:t></w:r>' +
          '</w:p>'
Set-RangeXml $combined $p45xml
